$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText {
    param($Address, $Text)
    $cell = $ws.Range($Address)
    $cell.Value = $Text
}

function Set-CellTextForceString {
    param($Address, $Text)
    $cell = $ws.Range($Address)
    # Temporarily force a text number format so a numeric-looking string
    # (e.g. "1.00") is not auto-converted into a number by Excel, then
    # restore the default "Normal" style so no extra style id is left on
    # the cell (matching the original unstyled inline-string cells).
    $cell.NumberFormat = "@"
    $cell.Value = $Text
    $cell.Style = "Normal"
}

Set-CellText "D2" '41.536.82'
Set-CellText "E2" '  +0.15%  '
Set-CellText "D3" '2.458.68'
Set-CellText "E3" '  +0.59%  '
Set-CellTextForceString "D4" '1.00'
Set-CellText "E4" '  -0.87%  '
Set-CellTextForceString "D5" '314.69'
Set-CellText "E5" '  +1.60%  '
Set-CellTextForceString "D6" '91.31'
Set-CellText "E6" '  +1.75%  '
Set-CellText "E7" '  +3.08%  '
Set-CellText "E8" '  -0.97%  '
Set-CellTextForceString "D9" '0.508'
Set-CellText "E9" '  +5.77%  '
Set-CellTextForceString "D10" '32.43'
Set-CellText "E10" '  +2.53%  '
Set-CellText "E11" '  +4.12%  '
Set-CellText "E12" '  +1.54%  '
Set-CellText "D13" '2.838.74'
Set-CellText "E13" '  +0.80%  '
Set-CellTextForceString "D14" '6.83'
Set-CellText "E14" '  +2.53%  '
Set-CellTextForceString "D15" '15.84'
Set-CellText "E15" '  +5.48%  '
Set-CellText "D16" '2.462.86'
Set-CellText "E16" '  +2.35%  '
Set-CellTextForceString "D17" '0.770'
Set-CellText "E17" '  +2.64%  '
Set-CellText "D18" '41.547.47'
Set-CellText "E18" '  +1.00%  '
Set-CellTextForceString "D19" '6.46'
Set-CellText "E19" '  +5.30%  '
Set-CellText "D20" '0.0₃0935'
Set-CellText "E20" '  +3.92%  '
Set-CellText "E21" '  +2.97%  '
Set-CellTextForceString "D22" '11.28'
Set-CellText "E22" '  +5.21%  '
Set-CellTextForceString "D23" '236.80'
Set-CellText "E23" '  +2.82%  '
Set-CellText "E24" '  +2.24%  '
Set-CellText "B25" 'Dai'
Set-CellText "C25" 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-CellTextForceString "D25" '1.00'
Set-CellText "E25" '  -0.21%  '
Set-CellText "B26" 'ImmutableX'
Set-CellText "C26" 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-CellTextForceString "D26" '1.90'
Set-CellText "E26" '  +3.27%  '
Set-CellTextForceString "D27" '24.25'
Set-CellText "E27" '  +3.13%  '
Set-CellTextForceString "D28" '2.25'
Set-CellText "E28" '  +1.91%  '
Set-CellTextForceString "D29" '9.64'
Set-CellText "E29" '  +2.30%  '
Set-CellTextForceString "D30" '34.98'
Set-CellText "E30" '  +0.02%  '
Set-CellTextForceString "D31" '155.54'
Set-CellText "E31" '  +2.99%  '
Set-CellTextForceString "D32" '5.42'
Set-CellText "E32" '  +3.75%  '
Set-CellText "E33" '  +2.34%  '
Set-CellTextForceString "D34" '0.0757'
Set-CellText "E34" '  +1.95%  '
Set-CellTextForceString "D35" '17.35'
Set-CellText "E35" '  -0.04%  '
Set-CellTextForceString "D36" '2.40'
Set-CellTextForceString "D37" '2.87'
Set-CellText "E37" '  -0.50%  '
Set-CellText "E38" '  +3.45%  '
Set-CellText "E39" '  +3.61%  '
Set-CellTextForceString "D40" '1.77'
Set-CellText "E40" '  -0.60%  '
Set-CellTextForceString "D41" '3.93'
Set-CellText "E41" '  -1.34%  '
Set-CellText "E42" '  -1.23%  '
Set-CellText "D43" '1.967.19'
Set-CellText "E43" '  +2.54%  '
Set-CellText "E44" '  +3.18%  '
Set-CellTextForceString "D45" '18.50'
Set-CellText "E45" '  -8.49%  '
Set-CellTextForceString "D46" '2.90'
Set-CellText "E46" '  +1.46%  '
Set-CellTextForceString "D47" '8.94'
Set-CellText "E47" '  +4.57%  '
Set-CellText "D48" '2.696.66'
Set-CellText "E48" '  +0.84%  '
Set-CellTextForceString "D49" '96.25'
Set-CellText "E49" '  +3.00%  '
Set-CellTextForceString "D50" '66.20'
Set-CellText "E50" '  +1.80%  '
Set-CellTextForceString "D51" '0.171'
Set-CellText "E51" '  -0.32%  '
